$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "TestA"
$ws.Range("B11").Value = 5000
$ws.Range("C11").Value = "New Business"
$ws.Range("F11").Value = "Follow up"
$ws.Range("G11").Value = "qualification"
$ws.Range("I11").Value = "Regina"

[void]$ws.Range("I11").Select()
